# Applies the cryptos list price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number must be forced to remain text
# (matching the original inlineStr/shared-string cell type) by using a leading
# apostrophe via .Formula, which Excel treats as a text quote-prefix instead of
# re-typing the cell as numeric.

$ws.Range("D2").Value = '29.435.83'
$ws.Range("E2").Value = '  +3.21%  '
$ws.Range("D3").Value = '1.604.98'
$ws.Range("E3").Value = '  +2.87%  '
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").Formula = "'212.28"
$ws.Range("E5").Value = '  +1.01%  '
$ws.Range("E6").Value = '  +6.90%  '
$ws.Range("D7").Formula = "'0.998"
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").Formula = "'26.77"
$ws.Range("E8").Value = '  +7.84%  '
$ws.Range("D9").Formula = "'43.48"
$ws.Range("E9").Value = '  -1.03%  '
$ws.Range("E10").Value = '  +2.86%  '
$ws.Range("E11").Value = '  +2.65%  '
$ws.Range("E12").Value = '  +1.50%  '
$ws.Range("D13").Value = '1.835.43'
$ws.Range("E13").Value = '  +2.88%  '
$ws.Range("D14").Value = '1.578.98'
$ws.Range("E14").Value = '  +1.45%  '
$ws.Range("D15").Value = '29.446.63'
$ws.Range("E15").Value = '  +3.14%  '
$ws.Range("E16").Value = '  +4.28%  '
$ws.Range("D17").Formula = "'3.70"
$ws.Range("E17").Value = '  +2.16%  '
$ws.Range("D18").Formula = "'63.07"
$ws.Range("E18").Value = '  +3.17%  '
$ws.Range("D19").Formula = "'241.84"
$ws.Range("E19").Value = '  +5.40%  '
$ws.Range("E20").Value = '  +4.18%  '
$ws.Range("D21").Value = '0.0₃0689'
$ws.Range("E21").Value = '  +2.02%  '
$ws.Range("D22").Formula = "'0.998"
$ws.Range("E22").Value = '  -0.22%  '
$ws.Range("E23").Value = '  +2.30%  '
$ws.Range("D24").Formula = "'9.18"
$ws.Range("E24").Value = '  +2.62%  '
$ws.Range("D26").Formula = "'154.50"
$ws.Range("E26").Value = '  +2.82%  '
$ws.Range("E27").Value = '  +5.37%  '
$ws.Range("D28").Formula = "'15.29"
$ws.Range("E28").Value = '  +3.65%  '
$ws.Range("D29").Formula = "'6.38"
$ws.Range("E29").Value = '  +2.76%  '
$ws.Range("E30").Value = '  -0.20%  '
$ws.Range("E32").Value = '  +0.14%  '
$ws.Range("E33").Value = '  +1.82%  '
$ws.Range("E34").Value = '  +4.83%  '
$ws.Range("D35").Value = '1.414.57'
$ws.Range("E35").Value = '  +2.00%  '
$ws.Range("D36").Formula = "'1.04"
$ws.Range("E36").Value = '  +0.40%  '
$ws.Range("E37").Value = '  +3.56%  '
$ws.Range("D38").Formula = "'2.83"
$ws.Range("E38").Value = '  +5.22%  '
$ws.Range("E39").Value = '  +0.26%  '
$ws.Range("D40").Formula = "'0.0166"
$ws.Range("E40").Value = '  +2.83%  '
$ws.Range("E41").Value = '  +3.97%  '
$ws.Range("E42").Value = '  +1.51%  '
$ws.Range("E43").Value = '  +6.31%  '
$ws.Range("D44").Formula = "'53.08"
$ws.Range("E44").Value = '  +23.08%  '
$ws.Range("D45").Formula = "'0.798"
$ws.Range("E45").Value = '  +3.47%  '
$ws.Range("E46").Value = '  -0.08%  '
$ws.Range("D47").Formula = "'65.75"
$ws.Range("E47").Value = '  +3.20%  '
$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").Formula = "'5.27"
$ws.Range("E48").Value = '  +0.93%  '
$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '1.745.51'
$ws.Range("E49").Value = '  +2.86%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Formula = "'86.84"
$ws.Range("E50").Value = '  +2.18%  '
$ws.Range("B51").Value = 'WEMIXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D51").Formula = "'0.847"
$ws.Range("E51").Value = '  -2.56%  '
